$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the two requisite lines: LOM3202 now comes first (row 24), LOM3206 second (row 25)
$ws.Range("B24").Value = "LOM3202 -  Circuitos Elétricos  (Requisito)`n"
$ws.Range("C24").Value = "LOM3202 -  Circuitos Elétricos  (Requisito)`n"
$ws.Range("B25").Value = "LOM3206 -  Eletrônica  (Indicação de Conjunto)`n"
$ws.Range("C25").Value = "LOM3206 -  Eletrônica  (Indicação de Conjunto)`n"
